# Auto-generated Excel COM-interop edit script
# Logs Week 16 values and performs season sim from Week 17
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 204
$ws.Range("F2").Value = 70
$ws.Range("G2").Value = 66
$ws.Range("H2").Value = 3
$ws.Range("J2").Value = 41
$ws.Range("L2").Value = 193
$ws.Range("M2").Value = 116
$ws.Range("O2").Value = 18
$ws.Range("Q2").Value = 440
$ws.Range("B3").Value = 13
$ws.Range("C3").Value = 131
$ws.Range("E3").Value = 26
$ws.Range("F3").Value = 83
$ws.Range("G3").Value = 22
$ws.Range("H3").Value = 20
$ws.Range("I3").Value = 34
$ws.Range("J3").Value = 54
$ws.Range("N3").Value = 15

$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 183
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 50
$ws.Range("G2").Value = 59
$ws.Range("H2").Value = 7
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 26
$ws.Range("L2").Value = 254
$ws.Range("M2").Value = 177
$ws.Range("O2").Value = 25
$ws.Range("P2").Value = 11
$ws.Range("Q2").Value = 458
$ws.Range("B3").Value = 11
$ws.Range("C3").Value = 169
$ws.Range("E3").Value = 29
$ws.Range("F3").Value = 110
$ws.Range("G3").Value = 35
$ws.Range("H3").Value = 25
$ws.Range("I3").Value = 58
$ws.Range("J3").Value = 49
$ws.Range("N3").Value = 6

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 82
$ws.Range("D2").Value = 49
$ws.Range("F2").Value = 180
$ws.Range("G2").Value = 171
$ws.Range("J2").Value = 70
$ws.Range("K2").Value = 64
$ws.Range("L2").Value = 43
$ws.Range("B3").Value = 54

$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 6
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 5

$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 14
$ws.Range("D2").Value = 14
$ws.Range("B3").Value = 14
$ws.Range("D3").Value = 6
$ws.Range("B4").Value = 4

# Shared long-text stat strings: append Week 16 numbers
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = "13 7 4 6 6 3 4 11 4 9 23 1 4 14 18 7 3 3 1 6 2 8 -6 -1 5 -3 9 -1 6 6 2 2 0 5 10 5 2 6 6 12 -2 9 -2 -3 1 7 9 7 6 6 13 3 3 0 8 9 27 1 -2 24 2 3 3 1 3 7 9 8 2 1 2 5 11 12 5 12 6 3 3 -4 7 2 2 7 4 2 13 10 6 4 1 3 0 1 4 4 1 4 14 3 6 1 6 7 18 -3 12 1 -2 11 1 1 -1 5 6 2 5 -2 9 11 23 14 7 3 -1 2 3 5 9 7 1 5 7 5 2 7 1 4 3 4 22 5 8 5 0 6 0 0 6 4 1 16 4 2 1 5 5 7 11 3 10 5 1 1 8 6 10 4 3 19 1 12 3 14 4 -2 20 4 7 3 10 4 3 21 4 3 5 2 2 3 0 3 8 18 4 2 5 4 2 1 4 6 8 5 4 5 3 5 6 2 3 10 2 2 3 3 0 5 2 12 2 1 10 2 4 6 0 11 10 17 -2 0 1 4 3 -8 7 4 2 5 11 31 0 3 4 3 0 5 -1 5 21 -2 15 5 25 1 7 5 3 1 2 6 1 23 1 2 3 4 -3 4 1 25 4 4 20 4 3 3 2 4 4 18 9 4 13 4 4 4 5 6 0 2 5 6 0 4 -3 7 5 10 6 2 -6 3 -1 7 2 1 14 1 3 24 -1 1 7 4 5 3 2 11 5 7 2 14 18 6 1 27 3 5 12 5 4 13 5 13 12 1 1 6 3 4 3 4 0 5 9 2 2 1 4 9 7 18 9 4 1 3 0 7 6 2 2 6 0 0 1 1 3 -1 3 5 -1 8 34 1 4 2 2 7 12 3 2 0 4 5 4 4 4 6 5 8 1 4 19 0 9 4 4 10 3 1 8 8 38 1 9 12 7 -3 4 3 15 4 2 3 4 4 0 8 4 -3 1 11 0 3 15 5 0 13 4 1 1 2 5 5 2 3 9 1 4 4 6 6 6 4 5 3 3 3 6 0"
$ws.Range("B3").Value = "10 11 2 18 13 28 8 4 19 3 9 -3 1 9 6 10 12 9 -1 5 8 25 7 19 9 0 23 4 11 4 9 1 5 26 91 6 10 9 14 38 5 5 9 19 7 19 28 41 27 0 3 4 1 5 12 5 7 24 19 16 12 15 -4 9 13 22 21 12 10 6 3 5 15 3 7 14 18 6 37 12 1 21 8 7 6 4 14 14 9 18 19 17 9 9 15 -7 6 -5 3 10 8 8 2 7 2 0 1 15 17 7 11 8 13 5 53 14 20 23 5 6 4 16 2 1 6 -3 24 25 6 24 13 8 -3 20 27 25 15 20 7 5 6 17 10 8 7 6 7 16 0 15 13 17 19 43 15 2 18 27 7 23 19 21 -6 19 4 28 16 18 3 4 6 24 36 8 6 4 5 9 11 12 14 2 8 4 16 14 33 -2 4 14 6 14 14 -1 23 14 0 19 2 0 13 4 21 6 10 8 7 14 11 16 9 36 12 10 12 5 25 13 16 12 22 -2 2 28 9 7 7 2 1 5 -1 19 28 16 45 6 5 8 17 21 20 5 9 7 9 7 34 19 17 6 2 3 23 13 46 1 9 1 6 39 4 3 16 10 5"
$ws.Range("C2").Value = "12 4 3 -2 3 10 11 14 1 8 7 12 3 2 11 3 1 4 3 -1 0 0 -1 3 8 5 -4 -1 3 3 0 1 2 5 2 7 4 5 2 2 10 2 3 2 7 4 21 2 0 1 3 9 3 0 8 -1 8 1 1 3 5 9 1 4 4 8 11 7 5 -1 6 5 7 6 4 9 5 8 1 2 10 3 0 0 13 2 -1 5 4 8 3 1 5 4 1 5 -6 5 4 5 4 5 3 12 2 15 24 17 1 1 2 9 3 10 10 1 8 1 2 4 1 12 9 16 3 3 5 4 5 1 5 6 0 26 3 1 5 0 3 3 -4 1 2 9 -5 8 2 3 11 1 8 6 11 3 5 -1 2 5 4 6 2 4 0 9 4 3 6 0 2 3 2 2 20 5 6 2 7 0 1 9 2 2 0 1 -1 3 1 5 1 2 3 12 9 4 3 9 4 8 -3 6 11 11 4 1 5 1 2 7 7 5 1 3 1 1 4 1 1 1 -2 5 4 2 2 -2 5 0 2 7 8 3 1 3 7 3 0 2 2 0 9 -3 0 1 1 5 -1 1 14 8 10 2 1 0 16 1 3 9 11 2 15 7 1 1 4 20 -4 8 9 1 2 6 -1 4 2 14 9 1 -2 4 3 1 9 1 -1 28 3 1 17 7 3 1 5 0 2 14 2 2 1 1 32 6 -1 -3 7 -4 3 4 13 -1 1 4 3 5 0 1 -1 2 1 3 17 9 1 0 1 15 2 2 3 2 3 8 17 0 2 6 5 2 4 1 3 3 5 4 -1 -1 8 -6 2 30 0 1 0 13 1 5 5 3 4 2 1 6 6 -2 -1 -4 10 3 7 1 6 5 0 5 0 2 2 4 4 2 7"
$ws.Range("C3").Value = "16 1 11 12 7 3 7 9 9 5 18 5 12 6 9 7 6 8 7 5 10 1 8 7 6 11 17 5 8 5 6 2 40 11 9 4 6 6 6 0 15 8 8 8 44 7 8 4 12 19 9 8 6 9 6 15 15 13 11 7 5 22 11 9 36 1 9 2 19 32 4 6 4 9 2 15 7 6 3 14 7 10 5 12 11 44 6 13 14 5 6 12 11 -1 8 5 7 0 7 5 24 2 8 9 4 7 25 9 22 16 10 2 8 8 11 3 9 23 -4 10 12 14 18 10 3 6 -4 5 4 4 9 15 0 8 5 4 8 10 2 27 10 12 4 43 -1 16 2 8 13 0 6 18 18 8 10 10 6 20 7 6 5 4 21 7 9 10 24 5 -1 29 3 18 7 5 8 -1 5 7 7 2 10 5 9 18 4 34 0 8 6 11 11 7 9 17 5 19 12 9 17 23 9 23 2 12 15 3 1 4 11 9 16 15 19 49 8 8 -2 7 3 8 18 3 14 4 5 11 -2 7 5 64 4 13 15 26 10 7 4 10 6 15 7 5 -9 0 13 9 14 3 5 7 6 18 19 5 10 9 26 21 -6 17 2 5 2 10 1 20 6 5 22 4 1 20 7 18 19 6 14 9 0 10 17 17 1 18 18 5 4 10 8 8 7 3 10 7 5 10 29 4 19 5 18 14 19 11 -2 15 4 14 7 11 19 7 6 46 5 1 7 5 9 7 5 12 29 5 6 1 13 10 3 -2 -4 8 8 8 -2 7 6 8 2 9 11 5 2 0 4 5 9 6 4 8 2"

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B4").Value = "73 23 62 70 69 67 68 66 58 57 57 65 73 73 60 63 65 65 62 66 65 65 61 67 64 62 65 61"
$ws.Range("B5").Value = "25 4 11 28 22 32 28 39 13 31 26 30 41 29 19 29 39 19 22 20 22 79 28 32 47 28 17 21"
$ws.Range("B6").Value = "18 22 20 14 44 27 19 26 23 23 24 17 19 18 14 20 13 15 26 21 5 13 31 23 17 20 10"
$ws.Range("D3").Value = "41 48 46 54 50 54 32 68 49 41 40 38 39 51 41 64 49 60 16 46 47 39 39 40 55 46 45 52 47 44 44 42 54 53 47 44 53 41 43 42 46 42 38 41 34 39 41 39 51"
$ws.Range("D4").Value = "0 7 0 7 0 13 0 0 5 8 0 0 0 0 0 0 22 8 0 9 12 0 4 0 5 0 0 15 11 0 10 0 8 19 0 8 0 14 7 9 0 0 0 0 0 0 3 13 0"
$ws.Range("D5").Value = "0 2 11 0 7 -1 11 2 0 0 0 -4 -3 0 5 0 0 0 22 0 0 0 0 8 6 0 0 0 12 0 10 0 9 12 10 0 8 0 13 4 0 2 3 9 0 0 11 0 1 10 -4 0 0 14 0 39 0 1"
